# Daily Update 2020 Feb 3
# Insert a new row for "San Benito, CA" (United States) into the lat/long
# table on Sheet1, just above the existing "Santa Clara, CA" row, pushing
# the remaining United States / Vietnam rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data table starts at row 3 and "Santa Clara, CA" currently lives at
# row 74. Insert a fresh row there so the new record lands above it.
$ws.Rows("74:74").Insert()

# Fill in the new row: State, Country, lat, long
$ws.Cells.Item(74, 1).Value = "San Benito, CA"
$ws.Cells.Item(74, 2).Value = "United States"
$ws.Cells.Item(74, 3).Value = 36.5761
$ws.Cells.Item(74, 4).Value = -120.9876

# Reflect the cursor/selection position left behind after the edit.
$ws.Range("B81").Select()
